$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range("G19").Value = 0.6924242424242424

# Row 21
$ws.Range("F21").Value = 0.5681818181818181
$ws.Range("G21").Value = 0.6924242424242424

# Row 22
$ws.Range("F22").Value = 0.8789473684210526

# Row 25
$ws.Range("F25").Value = 0.8578947368421053

# Row 26
$ws.Range("F26").Value = 0.7820895522388059

# Row 27
$ws.Range("F27").Value = 0.7432835820895523

# Row 28
$ws.Range("F28").Value = 0.7373134328358208

# Row 29
$ws.Range("F29").Value = 0.7343283582089553

# Row 31
$ws.Range("G31").Value = 0.88

# Row 32
$ws.Range("F32").Value = 0.76

# Row 33
$ws.Range("F33").Value = 0.744

# Row 41
$ws.Range("F41").Value = 0.8444444444444444
